$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update quarter header labels (row 8, 27, 46, 58, 77): drop oldest quarter, shift left, add new quarter ---
$headerRows = @(8, 27, 46, 58, 77)
$quarterLabels = @(
    "فصل چهارم منتهی به 1399/09",
    "فصل اول منتهی به 1399/12",
    "فصل دوم منتهی به 1400/03",
    "فصل سوم منتهی به 1400/06",
    "فصل چهارم منتهی به 1400/09",
    "فصل اول منتهی به 1400/12",
    "فصل دوم منتهی به 1401/03",
    "فصل سوم منتهی به 1401/06",
    "فصل چهارم منتهی به 1401/09",
    "فصل اول منتهی به 1401/12"
)
foreach ($r in $headerRows) {
    for ($i = 0; $i -lt 10; $i++) {
        $ws.Cells.Item($r, 5 + $i).Value = $quarterLabels[$i]
    }
}

# --- Update data rows: shift quarter data left by one column, append newest quarter value ---
# Row 11
$ws.Cells.Item(11, 5).Value = 0
$ws.Cells.Item(11, 6).Value = 0
$ws.Cells.Item(11, 7).Value = 0
$ws.Cells.Item(11, 8).Value = 0
$ws.Cells.Item(11, 9).Value = 0
$ws.Cells.Item(11, 10).Value = "-"
$ws.Cells.Item(11, 11).Value = "-"
$ws.Cells.Item(11, 12).Value = "-"
$ws.Cells.Item(11, 13).Value = "-"
$ws.Cells.Item(11, 14).Value = "-"

# Row 12
$ws.Cells.Item(12, 5).Value = 289497
$ws.Cells.Item(12, 6).Value = 320096
$ws.Cells.Item(12, 7).Value = 315433
$ws.Cells.Item(12, 8).Value = 223187
$ws.Cells.Item(12, 9).Value = 228338
$ws.Cells.Item(12, 10).Value = 218375
$ws.Cells.Item(12, 11).Value = 224591
$ws.Cells.Item(12, 12).Value = 243920
$ws.Cells.Item(12, 13).Value = 225074
$ws.Cells.Item(12, 14).Value = 210211

# Row 13
$ws.Cells.Item(13, 5).Value = 325268
$ws.Cells.Item(13, 6).Value = 287839
$ws.Cells.Item(13, 7).Value = 247378
$ws.Cells.Item(13, 8).Value = 182041
$ws.Cells.Item(13, 9).Value = 235990
$ws.Cells.Item(13, 10).Value = 232356
$ws.Cells.Item(13, 11).Value = 244703
$ws.Cells.Item(13, 12).Value = 249176
$ws.Cells.Item(13, 13).Value = 293430
$ws.Cells.Item(13, 14).Value = 221998

# Row 14
$ws.Cells.Item(14, 5).Value = 614765
$ws.Cells.Item(14, 6).Value = 607935
$ws.Cells.Item(14, 7).Value = 562811
$ws.Cells.Item(14, 8).Value = 405228
$ws.Cells.Item(14, 9).Value = 464328
$ws.Cells.Item(14, 10).Value = 450731
$ws.Cells.Item(14, 11).Value = 469294
$ws.Cells.Item(14, 12).Value = 493096
$ws.Cells.Item(14, 13).Value = 518504
$ws.Cells.Item(14, 14).Value = 432209

# Row 16
$ws.Cells.Item(16, 5).Value = 0
$ws.Cells.Item(16, 6).Value = 0
$ws.Cells.Item(16, 7).Value = 6799
$ws.Cells.Item(16, 8).Value = 2192
$ws.Cells.Item(16, 9).Value = 0
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 11).Value = 4000
$ws.Cells.Item(16, 12).Value = 5500
$ws.Cells.Item(16, 13).Value = 3900
$ws.Cells.Item(16, 14).Value = 3000

# Row 17
$ws.Cells.Item(17, 5).Value = 8300
$ws.Cells.Item(17, 6).Value = 22442
$ws.Cells.Item(17, 7).Value = 75564
$ws.Cells.Item(17, 8).Value = 12268
$ws.Cells.Item(17, 9).Value = 29099
$ws.Cells.Item(17, 10).Value = 18990
$ws.Cells.Item(17, 11).Value = 31230
$ws.Cells.Item(17, 12).Value = 53548
$ws.Cells.Item(17, 13).Value = 53021
$ws.Cells.Item(17, 14).Value = 9803

# Row 18
$ws.Cells.Item(18, 5).Value = 8300
$ws.Cells.Item(18, 6).Value = 22442
$ws.Cells.Item(18, 7).Value = 82363
$ws.Cells.Item(18, 8).Value = 14460
$ws.Cells.Item(18, 9).Value = 29099
$ws.Cells.Item(18, 10).Value = 18990
$ws.Cells.Item(18, 11).Value = 35230
$ws.Cells.Item(18, 12).Value = 59048
$ws.Cells.Item(18, 13).Value = 56921
$ws.Cells.Item(18, 14).Value = 12803

# Row 20
$ws.Cells.Item(20, 5).Value = 0
$ws.Cells.Item(20, 6).Value = 0
$ws.Cells.Item(20, 7).Value = 0
$ws.Cells.Item(20, 8).Value = 0
$ws.Cells.Item(20, 9).Value = 0
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(20, 11).Value = 0
$ws.Cells.Item(20, 12).Value = 0
$ws.Cells.Item(20, 13).Value = 0
$ws.Cells.Item(20, 14).Value = 0

# Row 21
$ws.Cells.Item(21, 5).Value = 0
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(21, 7).Value = 0
$ws.Cells.Item(21, 8).Value = "-"
$ws.Cells.Item(21, 9).Value = "-"
$ws.Cells.Item(21, 10).Value = "-"
$ws.Cells.Item(21, 11).Value = "-"
$ws.Cells.Item(21, 12).Value = "-"
$ws.Cells.Item(21, 13).Value = "-"
$ws.Cells.Item(21, 14).Value = "-"

# Row 22
$ws.Cells.Item(22, 5).Value = 0
$ws.Cells.Item(22, 6).Value = 0
$ws.Cells.Item(22, 7).Value = 0
$ws.Cells.Item(22, 8).Value = 0
$ws.Cells.Item(22, 9).Value = 0
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 0
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = 0
$ws.Cells.Item(22, 14).Value = 0

# Row 23
$ws.Cells.Item(23, 5).Value = 623065
$ws.Cells.Item(23, 6).Value = 630377
$ws.Cells.Item(23, 7).Value = 645174
$ws.Cells.Item(23, 8).Value = 419688
$ws.Cells.Item(23, 9).Value = 493427
$ws.Cells.Item(23, 10).Value = 469721
$ws.Cells.Item(23, 11).Value = 504524
$ws.Cells.Item(23, 12).Value = 552144
$ws.Cells.Item(23, 13).Value = 575425
$ws.Cells.Item(23, 14).Value = 445012

# Row 30
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(30, 6).Value = 0
$ws.Cells.Item(30, 7).Value = 0
$ws.Cells.Item(30, 8).Value = 0
$ws.Cells.Item(30, 9).Value = 0
$ws.Cells.Item(30, 10).Value = "-"
$ws.Cells.Item(30, 11).Value = "-"
$ws.Cells.Item(30, 12).Value = "-"
$ws.Cells.Item(30, 13).Value = "-"
$ws.Cells.Item(30, 14).Value = "-"

# Row 31
$ws.Cells.Item(31, 5).Value = 687352
$ws.Cells.Item(31, 6).Value = 708097
$ws.Cells.Item(31, 7).Value = 890038
$ws.Cells.Item(31, 8).Value = 1104535
$ws.Cells.Item(31, 9).Value = 998285
$ws.Cells.Item(31, 10).Value = 1105845
$ws.Cells.Item(31, 11).Value = 1166556
$ws.Cells.Item(31, 12).Value = 1794839
$ws.Cells.Item(31, 13).Value = 1543689
$ws.Cells.Item(31, 14).Value = 1357567

# Row 32
$ws.Cells.Item(32, 5).Value = 1004380
$ws.Cells.Item(32, 6).Value = 893535
$ws.Cells.Item(32, 7).Value = 984374
$ws.Cells.Item(32, 8).Value = 1046060
$ws.Cells.Item(32, 9).Value = 1341714
$ws.Cells.Item(32, 10).Value = 1378252
$ws.Cells.Item(32, 11).Value = 1550407
$ws.Cells.Item(32, 12).Value = 1905806
$ws.Cells.Item(32, 13).Value = 2104938
$ws.Cells.Item(32, 14).Value = 1642087

# Row 33
$ws.Cells.Item(33, 5).Value = 1691732
$ws.Cells.Item(33, 6).Value = 1601632
$ws.Cells.Item(33, 7).Value = 1874412
$ws.Cells.Item(33, 8).Value = 2150595
$ws.Cells.Item(33, 9).Value = 2339999
$ws.Cells.Item(33, 10).Value = 2484097
$ws.Cells.Item(33, 11).Value = 2716963
$ws.Cells.Item(33, 12).Value = 3700645
$ws.Cells.Item(33, 13).Value = 3648627
$ws.Cells.Item(33, 14).Value = 2999654

# Row 35
$ws.Cells.Item(35, 5).Value = 0
$ws.Cells.Item(35, 6).Value = 0
$ws.Cells.Item(35, 7).Value = 28237
$ws.Cells.Item(35, 8).Value = 11155
$ws.Cells.Item(35, 9).Value = 0
$ws.Cells.Item(35, 10).Value = 0
$ws.Cells.Item(35, 11).Value = 32576
$ws.Cells.Item(35, 12).Value = 42722
$ws.Cells.Item(35, 13).Value = 29978
$ws.Cells.Item(35, 14).Value = 23174

# Row 36
$ws.Cells.Item(36, 5).Value = 51966
$ws.Cells.Item(36, 6).Value = 125403
$ws.Cells.Item(36, 7).Value = 438125
$ws.Cells.Item(36, 8).Value = 99611
$ws.Cells.Item(36, 9).Value = 238920
$ws.Cells.Item(36, 10).Value = 148598
$ws.Cells.Item(36, 11).Value = 264296
$ws.Cells.Item(36, 12).Value = 507140
$ws.Cells.Item(36, 13).Value = 509973
$ws.Cells.Item(36, 14).Value = 116980

# Row 37
$ws.Cells.Item(37, 5).Value = 51966
$ws.Cells.Item(37, 6).Value = 125403
$ws.Cells.Item(37, 7).Value = 466362
$ws.Cells.Item(37, 8).Value = 110766
$ws.Cells.Item(37, 9).Value = 238920
$ws.Cells.Item(37, 10).Value = 148598
$ws.Cells.Item(37, 11).Value = 296872
$ws.Cells.Item(37, 12).Value = 549862
$ws.Cells.Item(37, 13).Value = 539951
$ws.Cells.Item(37, 14).Value = 140154

# Row 39
$ws.Cells.Item(39, 5).Value = 0
$ws.Cells.Item(39, 6).Value = 0
$ws.Cells.Item(39, 7).Value = 0
$ws.Cells.Item(39, 8).Value = 0
$ws.Cells.Item(39, 9).Value = 0
$ws.Cells.Item(39, 10).Value = 0
$ws.Cells.Item(39, 11).Value = 0
$ws.Cells.Item(39, 12).Value = 0
$ws.Cells.Item(39, 13).Value = 0
$ws.Cells.Item(39, 14).Value = 0

# Row 40
$ws.Cells.Item(40, 5).Value = 0
$ws.Cells.Item(40, 6).Value = 0
$ws.Cells.Item(40, 7).Value = 0
$ws.Cells.Item(40, 8).Value = "-"
$ws.Cells.Item(40, 9).Value = "-"
$ws.Cells.Item(40, 10).Value = "-"
$ws.Cells.Item(40, 11).Value = "-"
$ws.Cells.Item(40, 12).Value = "-"
$ws.Cells.Item(40, 13).Value = "-"
$ws.Cells.Item(40, 14).Value = "-"

# Row 41
$ws.Cells.Item(41, 5).Value = -87158
$ws.Cells.Item(41, 6).Value = 0
$ws.Cells.Item(41, 7).Value = -3552
$ws.Cells.Item(41, 8).Value = 0
$ws.Cells.Item(41, 9).Value = 0
$ws.Cells.Item(41, 10).Value = 0
$ws.Cells.Item(41, 11).Value = 0
$ws.Cells.Item(41, 12).Value = 0
$ws.Cells.Item(41, 13).Value = 0
$ws.Cells.Item(41, 14).Value = 0

# Row 42
$ws.Cells.Item(42, 5).Value = 1656540
$ws.Cells.Item(42, 6).Value = 1727035
$ws.Cells.Item(42, 7).Value = 2337222
$ws.Cells.Item(42, 8).Value = 2261361
$ws.Cells.Item(42, 9).Value = 2578919
$ws.Cells.Item(42, 10).Value = 2632695
$ws.Cells.Item(42, 11).Value = 3013835
$ws.Cells.Item(42, 12).Value = 4250507
$ws.Cells.Item(42, 13).Value = 4188578
$ws.Cells.Item(42, 14).Value = 3139808

# Row 49
$ws.Cells.Item(49, 5).Value = "-"
$ws.Cells.Item(49, 6).Value = "-"
$ws.Cells.Item(49, 7).Value = "-"
$ws.Cells.Item(49, 8).Value = "-"
$ws.Cells.Item(49, 9).Value = "-"
$ws.Cells.Item(49, 10).Value = "-"
$ws.Cells.Item(49, 11).Value = "-"
$ws.Cells.Item(49, 12).Value = "-"
$ws.Cells.Item(49, 13).Value = "-"
$ws.Cells.Item(49, 14).Value = "-"

# Row 50
$ws.Cells.Item(50, 5).Value = 2374297
$ws.Cells.Item(50, 6).Value = 2212139
$ws.Cells.Item(50, 7).Value = 2821639
$ws.Cells.Item(50, 8).Value = 4948922
$ws.Cells.Item(50, 9).Value = 4371962
$ws.Cells.Item(50, 10).Value = 5063973
$ws.Cells.Item(50, 11).Value = 5129967
$ws.Cells.Item(50, 12).Value = 5921274
$ws.Cells.Item(50, 13).Value = 6858584
$ws.Cells.Item(50, 14).Value = 6458117

# Row 51
$ws.Cells.Item(51, 5).Value = 3087854
$ws.Cells.Item(51, 6).Value = 3104287
$ws.Cells.Item(51, 7).Value = 3979230
$ws.Cells.Item(51, 8).Value = 5746288
$ws.Cells.Item(51, 9).Value = 5685470
$ws.Cells.Item(51, 10).Value = 5931639
$ws.Cells.Item(51, 11).Value = 6138987
$ws.Cells.Item(51, 12).Value = 6656888
$ws.Cells.Item(51, 13).Value = 7173561
$ws.Cells.Item(51, 14).Value = 7396852

# Row 53
$ws.Cells.Item(53, 5).Value = "-"
$ws.Cells.Item(53, 6).Value = "-"
$ws.Cells.Item(53, 7).Value = 4153111
$ws.Cells.Item(53, 8).Value = 5088960
$ws.Cells.Item(53, 9).Value = "-"
$ws.Cells.Item(53, 10).Value = "-"
$ws.Cells.Item(53, 11).Value = 8144000
$ws.Cells.Item(53, 12).Value = 7926105
$ws.Cells.Item(53, 13).Value = 7686667
$ws.Cells.Item(53, 14).Value = 7724667

# Row 54
$ws.Cells.Item(54, 5).Value = 6260964
$ws.Cells.Item(54, 6).Value = 5587871
$ws.Cells.Item(54, 7).Value = 5798065
$ws.Cells.Item(54, 8).Value = 8119579
$ws.Cells.Item(54, 9).Value = 8210591
$ws.Cells.Item(54, 10).Value = 7825066
$ws.Cells.Item(54, 11).Value = 8221705
$ws.Cells.Item(54, 12).Value = 8866259
$ws.Cells.Item(54, 13).Value = 9618321
$ws.Cells.Item(54, 14).Value = 11932960

# Row 61
$ws.Cells.Item(61, 5).Value = 0
$ws.Cells.Item(61, 6).Value = 0
$ws.Cells.Item(61, 7).Value = 0
$ws.Cells.Item(61, 8).Value = 0
$ws.Cells.Item(61, 9).Value = 0
$ws.Cells.Item(61, 10).Value = "-"
$ws.Cells.Item(61, 11).Value = "-"
$ws.Cells.Item(61, 12).Value = "-"
$ws.Cells.Item(61, 13).Value = "-"
$ws.Cells.Item(61, 14).Value = "-"

# Row 62
$ws.Cells.Item(62, 5).Value = -289647
$ws.Cells.Item(62, 6).Value = -395997
$ws.Cells.Item(62, 7).Value = -362760
$ws.Cells.Item(62, 8).Value = -449150
$ws.Cells.Item(62, 9).Value = -469132
$ws.Cells.Item(62, 10).Value = -577640
$ws.Cells.Item(62, 11).Value = -385762
$ws.Cells.Item(62, 12).Value = -880935
$ws.Cells.Item(62, 13).Value = -498214
$ws.Cells.Item(62, 14).Value = -677352

# Row 63
$ws.Cells.Item(63, 5).Value = -478085
$ws.Cells.Item(63, 6).Value = -499701
$ws.Cells.Item(63, 7).Value = -492848
$ws.Cells.Item(63, 8).Value = -515631
$ws.Cells.Item(63, 9).Value = -597873
$ws.Cells.Item(63, 10).Value = -719932
$ws.Cells.Item(63, 11).Value = -681678
$ws.Cells.Item(63, 12).Value = -1055606
$ws.Cells.Item(63, 13).Value = -937630
$ws.Cells.Item(63, 14).Value = -907573

# Row 64
$ws.Cells.Item(64, 5).Value = -767732
$ws.Cells.Item(64, 6).Value = -895698
$ws.Cells.Item(64, 7).Value = -855608
$ws.Cells.Item(64, 8).Value = -964781
$ws.Cells.Item(64, 9).Value = -1067005
$ws.Cells.Item(64, 10).Value = -1297572
$ws.Cells.Item(64, 11).Value = -1067440
$ws.Cells.Item(64, 12).Value = -1936541
$ws.Cells.Item(64, 13).Value = -1435844
$ws.Cells.Item(64, 14).Value = -1584925

# Row 66
$ws.Cells.Item(66, 5).Value = 0
$ws.Cells.Item(66, 6).Value = 0
$ws.Cells.Item(66, 7).Value = -7597
$ws.Cells.Item(66, 8).Value = -4042
$ws.Cells.Item(66, 9).Value = -555
$ws.Cells.Item(66, 10).Value = 0
$ws.Cells.Item(66, 11).Value = -7988
$ws.Cells.Item(66, 12).Value = -11141
$ws.Cells.Item(66, 13).Value = -11611
$ws.Cells.Item(66, 14).Value = -9082

# Row 67
$ws.Cells.Item(67, 5).Value = -14438
$ws.Cells.Item(67, 6).Value = -70131
$ws.Cells.Item(67, 7).Value = -204982
$ws.Cells.Item(67, 8).Value = -81020
$ws.Cells.Item(67, 9).Value = -82729
$ws.Cells.Item(67, 10).Value = -77620
$ws.Cells.Item(67, 11).Value = -120683
$ws.Cells.Item(67, 12).Value = -243567
$ws.Cells.Item(67, 13).Value = -202045
$ws.Cells.Item(67, 14).Value = -43046

# Row 68
$ws.Cells.Item(68, 5).Value = -14438
$ws.Cells.Item(68, 6).Value = -70131
$ws.Cells.Item(68, 7).Value = -212579
$ws.Cells.Item(68, 8).Value = -85062
$ws.Cells.Item(68, 9).Value = -83284
$ws.Cells.Item(68, 10).Value = -77620
$ws.Cells.Item(68, 11).Value = -128671
$ws.Cells.Item(68, 12).Value = -254708
$ws.Cells.Item(68, 13).Value = -213656
$ws.Cells.Item(68, 14).Value = -52128

# Row 70
$ws.Cells.Item(70, 5).Value = 0
$ws.Cells.Item(70, 6).Value = 0
$ws.Cells.Item(70, 7).Value = 0
$ws.Cells.Item(70, 8).Value = 0
$ws.Cells.Item(70, 9).Value = 0
$ws.Cells.Item(70, 10).Value = 0
$ws.Cells.Item(70, 11).Value = 0
$ws.Cells.Item(70, 12).Value = 0
$ws.Cells.Item(70, 13).Value = 0
$ws.Cells.Item(70, 14).Value = 0

# Row 71
$ws.Cells.Item(71, 5).Value = 0
$ws.Cells.Item(71, 6).Value = 0
$ws.Cells.Item(71, 7).Value = 0
$ws.Cells.Item(71, 8).Value = "-"
$ws.Cells.Item(71, 9).Value = "-"
$ws.Cells.Item(71, 10).Value = "-"
$ws.Cells.Item(71, 11).Value = "-"
$ws.Cells.Item(71, 12).Value = "-"
$ws.Cells.Item(71, 13).Value = "-"
$ws.Cells.Item(71, 14).Value = "-"

# Row 72
$ws.Cells.Item(72, 5).Value = 0
$ws.Cells.Item(72, 6).Value = 0
$ws.Cells.Item(72, 7).Value = 0
$ws.Cells.Item(72, 8).Value = 0
$ws.Cells.Item(72, 9).Value = 0
$ws.Cells.Item(72, 10).Value = 0
$ws.Cells.Item(72, 11).Value = 0
$ws.Cells.Item(72, 12).Value = 0
$ws.Cells.Item(72, 13).Value = 0
$ws.Cells.Item(72, 14).Value = 0

# Row 73
$ws.Cells.Item(73, 5).Value = -782170
$ws.Cells.Item(73, 6).Value = -965829
$ws.Cells.Item(73, 7).Value = -1068187
$ws.Cells.Item(73, 8).Value = -1049843
$ws.Cells.Item(73, 9).Value = -1150289
$ws.Cells.Item(73, 10).Value = -1375192
$ws.Cells.Item(73, 11).Value = -1196111
$ws.Cells.Item(73, 12).Value = -2191249
$ws.Cells.Item(73, 13).Value = -1649500
$ws.Cells.Item(73, 14).Value = -1637053

# Row 80
$ws.Cells.Item(80, 5).Value = 0
$ws.Cells.Item(80, 6).Value = 0
$ws.Cells.Item(80, 7).Value = 0
$ws.Cells.Item(80, 8).Value = 0
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 10).Value = "-"
$ws.Cells.Item(80, 11).Value = "-"
$ws.Cells.Item(80, 12).Value = "-"
$ws.Cells.Item(80, 13).Value = "-"
$ws.Cells.Item(80, 14).Value = "-"

# Row 81
$ws.Cells.Item(81, 5).Value = 352632
$ws.Cells.Item(81, 6).Value = 312100
$ws.Cells.Item(81, 7).Value = 527278
$ws.Cells.Item(81, 8).Value = 653343
$ws.Cells.Item(81, 9).Value = 529153
$ws.Cells.Item(81, 10).Value = 528205
$ws.Cells.Item(81, 11).Value = 780794
$ws.Cells.Item(81, 12).Value = 913904
$ws.Cells.Item(81, 13).Value = 1045475
$ws.Cells.Item(81, 14).Value = 680215

# Row 82
$ws.Cells.Item(82, 5).Value = 490212
$ws.Cells.Item(82, 6).Value = 393834
$ws.Cells.Item(82, 7).Value = 491526
$ws.Cells.Item(82, 8).Value = 528919
$ws.Cells.Item(82, 9).Value = 743841
$ws.Cells.Item(82, 10).Value = 658320
$ws.Cells.Item(82, 11).Value = 868729
$ws.Cells.Item(82, 12).Value = 850200
$ws.Cells.Item(82, 13).Value = 1167308
$ws.Cells.Item(82, 14).Value = 734514

# Row 83
$ws.Cells.Item(83, 5).Value = 842844
$ws.Cells.Item(83, 6).Value = 705934
$ws.Cells.Item(83, 7).Value = 1018804
$ws.Cells.Item(83, 8).Value = 1182262
$ws.Cells.Item(83, 9).Value = 1272994
$ws.Cells.Item(83, 10).Value = 1186525
$ws.Cells.Item(83, 11).Value = 1649523
$ws.Cells.Item(83, 12).Value = 1764104
$ws.Cells.Item(83, 13).Value = 2212783
$ws.Cells.Item(83, 14).Value = 1414729

# Row 85
$ws.Cells.Item(85, 5).Value = 0
$ws.Cells.Item(85, 6).Value = 0
$ws.Cells.Item(85, 7).Value = 20640
$ws.Cells.Item(85, 8).Value = 7113
$ws.Cells.Item(85, 9).Value = -555
$ws.Cells.Item(85, 10).Value = 0
$ws.Cells.Item(85, 11).Value = 24588
$ws.Cells.Item(85, 12).Value = 31581
$ws.Cells.Item(85, 13).Value = 18367
$ws.Cells.Item(85, 14).Value = 14092

# Row 86
$ws.Cells.Item(86, 5).Value = 37528
$ws.Cells.Item(86, 6).Value = 55272
$ws.Cells.Item(86, 7).Value = 233143
$ws.Cells.Item(86, 8).Value = 18591
$ws.Cells.Item(86, 9).Value = 156191
$ws.Cells.Item(86, 10).Value = 70978
$ws.Cells.Item(86, 11).Value = 143613
$ws.Cells.Item(86, 12).Value = 263573
$ws.Cells.Item(86, 13).Value = 307928
$ws.Cells.Item(86, 14).Value = 73934

# Row 87
$ws.Cells.Item(87, 5).Value = 37528
$ws.Cells.Item(87, 6).Value = 55272
$ws.Cells.Item(87, 7).Value = 253783
$ws.Cells.Item(87, 8).Value = 25704
$ws.Cells.Item(87, 9).Value = 155636
$ws.Cells.Item(87, 10).Value = 70978
$ws.Cells.Item(87, 11).Value = 168201
$ws.Cells.Item(87, 12).Value = 295154
$ws.Cells.Item(87, 13).Value = 326295
$ws.Cells.Item(87, 14).Value = 88026

# Row 89
$ws.Cells.Item(89, 5).Value = 0
$ws.Cells.Item(89, 6).Value = 0
$ws.Cells.Item(89, 7).Value = 0
$ws.Cells.Item(89, 8).Value = 0
$ws.Cells.Item(89, 9).Value = 0
$ws.Cells.Item(89, 10).Value = 0
$ws.Cells.Item(89, 11).Value = 0
$ws.Cells.Item(89, 12).Value = 0
$ws.Cells.Item(89, 13).Value = 0
$ws.Cells.Item(89, 14).Value = 0

# Row 90
$ws.Cells.Item(90, 5).Value = 880372
$ws.Cells.Item(90, 6).Value = 761206
$ws.Cells.Item(90, 7).Value = 1272587
$ws.Cells.Item(90, 8).Value = 1207966
$ws.Cells.Item(90, 9).Value = 1428630
$ws.Cells.Item(90, 10).Value = 1257503
$ws.Cells.Item(90, 11).Value = 1817724
$ws.Cells.Item(90, 12).Value = 2059258
$ws.Cells.Item(90, 13).Value = 2539078
$ws.Cells.Item(90, 14).Value = 1502755
